$d = $word.ActiveDocument

# Update the date heading in the first paragraph
$d.Content.Find.Execute("2025-10-22 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-10-23 Thursday", 2) | Out-Null

# Update each arithmetic-problem cell in the table, in row-major order
$t = $d.Tables.Item(1)
$values = @(
    "61+10=",
    "6+36=",
    "31-7=",
    "41-39=",
    "38+29=",
    "95-25=",
    "13+14=",
    "19-15=",
    "25+73=",
    "2+4=",
    "58+37=",
    "91-54=",
    "30-21=",
    "8+60=",
    "79-20=",
    "90-86=",
    "58-16=",
    "73-0=",
    "14+33=",
    "36+43=",
    "86-6=",
    "87-4=",
    "29+17=",
    "17+52=",
    "60-1=",
    "6+12=",
    "19-5=",
    "35+11=",
    "13+36=",
    "47+35=",
    "40-36=",
    "46+6=",
    "35-7=",
    "49+30=",
    "19+43=",
    "11+78=",
    "80-64=",
    "94-15=",
    "5+45=",
    "25+74=",
    "3+88=",
    "58-34=",
    "28+56=",
    "89-76=",
    "96-70=",
    "94-48=",
    "59+27=",
    "37+27=",
    "0+72=",
    "0+81=",
    "79-5=",
    "23+68=",
    "3+14=",
    "74-26=",
    "56-0=",
    "4+68=",
    "7+47=",
    "11+52=",
    "25+53=",
    "14-9=",
    "50-5=",
    "18+54=",
    "46-38=",
    "83-67=",
    "6+67=",
    "10+50=",
    "17+32=",
    "11+42=",
    "62-30=",
    "5+15=",
    "16-0=",
    "90-50=",
    "26+23=",
    "49+17=",
    "56+18=",
    "61-58=",
    "85-13=",
    "69+25=",
    "80-60=",
    "22+44=",
    "27+54=",
    "30-15=",
    "40+44=",
    "82-33=",
    "90-16=",
    "52+3=",
    "53-52=",
    "82+10=",
    "97-60=",
    "55-34=",
    "92+4=",
    "21-8=",
    "95-75=",
    "41+5=",
    "59-44=",
    "90-67=",
    "78-19=",
    "58-25=",
    "91-14=",
    "36+0="
)

$rows = $t.Rows.Count
$cols = $t.Columns.Count
$idx = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $t.Cell($r, $c)
        $rng = $cell.Range
        $rng.End = $rng.End - 2
        $rng.Text = $values[$idx]
        $idx++
    }
}

Write-Host "Done. Updated" $idx "cells."
